# feat: create and implement logic for visual cues on login & signin
#
# Inserts two new rows (14-15) into the localisation table describing a new
# "UI Signup feedback" key, shifting the existing rows (old 14-29) down to
# 16-31. Mirrors the already-present "UI Login feedback" pattern (rows 12-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows, pushing everything from row 14 down by two.
$ws.Rows("14:15").Insert()

# Bring formatting (style s="4") along for the new rows by cloning row 16
# (the old row 14, now shifted down) onto the freshly inserted 14:15.
$ws.Range("A16:F16").Copy()
$ws.Range("A14:F15").PasteSpecial(-4122)

# Row 14: "valid signup" feedback key.
$ws.Range("A14").Value = "UI Signup feedback"
$ws.Range("B14").Value = "tmp valid signup"

# Row 15: "invalid signup" feedback key.
$ws.Range("A15").Value = "UI Signup feedback"
$ws.Range("B15").Value = "tmp invalid signup"

# French localisation text for the "valid" row.
$ws.Range("D14").Value = "Votre compte a été créé!"

# English/French localisation text for the "invalid" row.
$ws.Range("C15").Value = "Account could not be created."
$ws.Range("D15").Value = "Échec de la création du compte."

# English localisation text for the "valid" row (filled in last).
$ws.Range("C14").Value = "Account created!"

# The "invalid signup" row wraps onto two lines in Excel, same as the
# existing "invalid login" row (row 13).
$ws.Rows("15:15").RowHeight = 28.8

# Update the sheet's active selection to where the author ended up editing.
$ws.Range("C13").Select()
